# Added Registration Page and Refactored code.
#
# - Rename the third sheet ("Sheet3") to "Registration".
# - Populate it with a header row (first name / last name) that was
#   previously empty.
# - Make "Registration" the active sheet/tab (was "LoginDetails").
# - Move the selection on the Registration sheet to B4.

$wb = $excel.ActiveWorkbook

$wsRegistration = $wb.Worksheets.Item(3)
$wsRegistration.Name = "Registration"

$wsRegistration.Range("A1").Value = "Gowri"
$wsRegistration.Range("B1").Value = "Kumar"

# Activating this sheet makes it the workbook's active tab and clears
# "tabSelected" from whichever sheet used to carry it (LoginDetails).
$wsRegistration.Activate() | Out-Null
$wsRegistration.Range("B4").Select() | Out-Null
